# The upstream commit only normalizes the canonical XML serialization of
# word/document.xml and word/styles.xml (w:* attributes are re-emitted in
# alphabetical order by a repository tooling pass); no paragraph text,
# field code, run formatting, page geometry, style, or document-property
# value actually changes - before and after are identical once attribute
# order is ignored. Word's object model does not expose a way to control
# the raw attribute order that the XML writer uses, so there is no
# content edit to make here: simply touch and resave the document so the
# package is (re)written through the normal save path without altering
# any of its content.
$d = $word.ActiveDocument
$d.Save()
